# Updates the crypto price table (columns B-E) to match the latest
# coinranking.com snapshot. Column D holds prices formatted as plain text
# (e.g. '43.835.18'), so numeric-looking values are entered with a leading
# apostrophe to keep Excel from reinterpreting them as numbers, then the
# cell style is reset to Normal so no stray "Text" / quote-prefix formatting
# is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'43.835.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.36%  "

# Row 3
$ws.Range("D3").Value = "'2.337.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.34%  "

# Row 4
$ws.Range("D4").Value = "'1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'312.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "

# Row 6
$ws.Range("D6").Value = "'108.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.96%  "

# Row 8
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("D9").Value = "'0.620"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.30%  "

# Row 10
$ws.Range("E10").Value = "  +4.75%  "

# Row 11
$ws.Range("D11").Value = "'0.0920"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.80%  "

# Row 12
$ws.Range("D12").Value = "'8.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.55%  "

# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'1.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.88%  "

# Row 14
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "'0.108"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.17%  "

# Row 15
$ws.Range("D15").Value = "'15.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.08%  "

# Row 16
$ws.Range("D16").Value = "'2.692.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.34%  "

# Row 17
$ws.Range("D17").Value = "'2.336.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.08%  "

# Row 18
$ws.Range("D18").Value = "'43.781.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.47%  "

# Row 19
$ws.Range("E19").Value = "  +2.03%  "

# Row 20
$ws.Range("E20").Value = "  +1.13%  "

# Row 21
$ws.Range("D21").Value = "'13.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.50%  "

# Row 22
$ws.Range("D22").Value = "'74.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.58%  "

# Row 23
$ws.Range("D23").Value = "'3.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.57%  "

# Row 24
$ws.Range("D24").Value = "'267.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.13%  "

# Row 25
$ws.Range("E25").Value = "  +3.30%  "

# Row 26
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("D27").Value = "'7.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.27%  "

# Row 28
$ws.Range("D28").Value = "'11.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.01%  "

# Row 29
$ws.Range("E29").Value = "  -1.85%  "

# Row 30
$ws.Range("D30").Value = "'39.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.73%  "

# Row 31
$ws.Range("D31").Value = "'22.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.55%  "

# Row 32
$ws.Range("D32").Value = "'168.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.05%  "

# Row 33
$ws.Range("D33").Value = "'0.0890"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.15%  "

# Row 34
$ws.Range("D34").Value = "'2.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.94%  "

# Row 35
$ws.Range("E35").Value = "  +0.94%  "

# Row 36
$ws.Range("D36").Value = "'0.114"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.38%  "

# Row 37
$ws.Range("E37").Value = "  +4.05%  "

# Row 38
$ws.Range("D38").Value = "'0.0365"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.16%  "

# Row 39
$ws.Range("E39").Value = "  +8.74%  "

# Row 40
$ws.Range("D40").Value = "'3.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.29%  "

# Row 41
$ws.Range("E41").Value = "  +8.96%  "

# Row 42
$ws.Range("D42").Value = "'104.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.84%  "

# Row 43
$ws.Range("E43").Value = "  +3.25%  "

# Row 44
$ws.Range("D44").Value = "'13.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.82%  "

# Row 45
$ws.Range("D45").Value = "'71.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.87%  "

# Row 46
$ws.Range("E46").Value = "  -0.08%  "

# Row 47
$ws.Range("D47").Value = "'114.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.54%  "

# Row 48
$ws.Range("B48").Value = "TheGraph"
$ws.Range("C48").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D48").Value = "'0.221"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +17.18%  "

# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'1.660.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.62%  "

# Row 50
$ws.Range("D50").Value = "'76.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.61%  "

# Row 51
$ws.Range("D51").Value = "'8.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.96%  "
